# ---------------------------------------------------------------------------
# Applies the "Add files via upload" commit to youtube_sentiment_report.xlsx:
#   - Overall Summary (sheet 1): refresh the summary row (row 2) with new
#     totals, and append 7 rows (85-91) mirroring the newest "Video Details"
#     rows plus two "== Updated ... ===" timestamp markers.
#   - Video Details (sheet 2): append 3 new video rows (66-68).
#   - Positive Keywords (sheet 3): append 6 new keyword rows (104-109).
#   - Negative Keywords (sheet 4): append 12 new keyword rows (168-179).
#   - Negative Issues (sheet 7): append 11 new issue rows (127-137).
#   - Narrative Summary (sheet 9) and Political Issues Narrative (sheet 10):
#     replace the single cumulative-narrative cell with refreshed AI copy.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $text) {
    # Force literal-text storage so numeric/percent-looking strings ("100.0%",
    # "0.875", "-1.000", ...) are not auto-coerced into numbers by Excel.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

function Restore-RowStyle($templateRowRange, $targetRowRange) {
    # Re-apply the donor row's cell formatting (style index) on top of the
    # values we just wrote, without touching the values themselves.
    $templateRowRange.Copy()
    $targetRowRange.PasteSpecial(-4122)  # xlPasteFormats
}

# ===========================================================================
# Sheet 1: "Overall Summary"
# ===========================================================================
$ws1 = $wb.Worksheets.Item("Overall Summary")

# --- refresh cumulative summary row (row 2) ---
$ws1.Cells.Item(2,1).Value = 70
$ws1.Cells.Item(2,2).Value = 52
$ws1.Cells.Item(2,3).Value = 0
$ws1.Cells.Item(2,4).Value = 100
$ws1.Cells.Item(2,5).Value = 0
Set-TextCell $ws1 2 6 "0.0"
Set-TextCell $ws1 2 7 "100.0"
Set-TextCell $ws1 2 8 "0.0"
Set-TextCell $ws1 2 9 "-1.000"
$srcRow2 = $ws1.Range("F3:I3")
$dstRow2 = $ws1.Range("F2:I2")
Restore-RowStyle $srcRow2 $dstRow2

# --- append rows 85-91 ---
# row 85: blank separator (handled by final style pass below)
# row 86: timestamp marker formula
Set-TextCell $ws1 86 1 "placeholder"
$ws1.Cells.Item(86,1).Formula = "=== Updated 2025-10-16 20:13:50 ==="
# row 87: blank separator (handled by final style pass below)
# row 88: new video row
$ws1.Cells.Item(88,1).Value = "Venkata Rami Reddy Fires On Chandrababu | AP Gover"
$ws1.Cells.Item(88,2).Value = "Sakshi TV Live"
$ws1.Cells.Item(88,3).Value = 10
Set-TextCell $ws1 88 4 "100.0%"
Set-TextCell $ws1 88 5 "0.0%"
Set-TextCell $ws1 88 6 "0.0%"
$ws1.Cells.Item(88,7).Value = 10
$ws1.Cells.Item(88,8).Value = 4
# row 89: timestamp marker formula
Set-TextCell $ws1 89 1 "placeholder"
$ws1.Cells.Item(89,1).Formula = "=== Updated 2025-10-16 20:46:23 ==="
# row 90: new video row
$ws1.Cells.Item(90,1).Value = "Speaker Ayyannapatrudu Reads Jagan Letter In AP As"
$ws1.Cells.Item(90,2).Value = "Cloud Media"
$ws1.Cells.Item(90,3).Value = 15
Set-TextCell $ws1 90 4 "0.0%"
Set-TextCell $ws1 90 5 "100.0%"
Set-TextCell $ws1 90 6 "0.0%"
$ws1.Cells.Item(90,7).Value = 10
$ws1.Cells.Item(90,8).Value = 4
# row 91: new video row
$ws1.Cells.Item(91,1).Value = "Minister Atchannaidu Fires on YS Jagan | Onion Far"
Set-TextCell $ws1 91 2 "TV5 News "
$ws1.Cells.Item(91,3).Value = 10
Set-TextCell $ws1 91 4 "0.0%"
Set-TextCell $ws1 91 5 "100.0%"
Set-TextCell $ws1 91 6 "0.0%"
$ws1.Cells.Item(91,7).Value = 10
$ws1.Cells.Item(91,8).Value = 3

# Re-apply the plain data-row style (s=2, matches rows 3 / 84) over rows
# 85-91 now that every value is in place.
$ws1Src = $ws1.Range("A3:I3")
$ws1Dst85 = $ws1.Range("A85:I85")
Restore-RowStyle $ws1Src $ws1Dst85
$ws1Dst86 = $ws1.Range("A86:I86")
Restore-RowStyle $ws1Src $ws1Dst86
$ws1Dst87 = $ws1.Range("A87:I87")
Restore-RowStyle $ws1Src $ws1Dst87
$ws1Dst88 = $ws1.Range("A88:I88")
Restore-RowStyle $ws1Src $ws1Dst88
$ws1Dst89 = $ws1.Range("A89:I89")
Restore-RowStyle $ws1Src $ws1Dst89
$ws1Dst90 = $ws1.Range("A90:I90")
Restore-RowStyle $ws1Src $ws1Dst90
$ws1Dst91 = $ws1.Range("A91:I91")
Restore-RowStyle $ws1Src $ws1Dst91

# ===========================================================================
# Sheet 2: "Video Details"
# ===========================================================================
$ws2 = $wb.Worksheets.Item("Video Details")

$ws2.Cells.Item(66,1).Value = "Venkata Rami Reddy Fires On Chandrababu | AP Gover"
$ws2.Cells.Item(66,2).Value = "Sakshi TV Live"
$ws2.Cells.Item(66,3).Value = 10
Set-TextCell $ws2 66 4 "100.0%"
Set-TextCell $ws2 66 5 "0.0%"
Set-TextCell $ws2 66 6 "0.0%"
$ws2.Cells.Item(66,7).Value = 10
$ws2.Cells.Item(66,8).Value = 4

$ws2.Cells.Item(67,1).Value = "Speaker Ayyannapatrudu Reads Jagan Letter In AP As"
$ws2.Cells.Item(67,2).Value = "Cloud Media"
$ws2.Cells.Item(67,3).Value = 15
Set-TextCell $ws2 67 4 "0.0%"
Set-TextCell $ws2 67 5 "100.0%"
Set-TextCell $ws2 67 6 "0.0%"
$ws2.Cells.Item(67,7).Value = 10
$ws2.Cells.Item(67,8).Value = 4

$ws2.Cells.Item(68,1).Value = "Minister Atchannaidu Fires on YS Jagan | Onion Far"
Set-TextCell $ws2 68 2 "TV5 News "
$ws2.Cells.Item(68,3).Value = 10
Set-TextCell $ws2 68 4 "0.0%"
Set-TextCell $ws2 68 5 "100.0%"
Set-TextCell $ws2 68 6 "0.0%"
$ws2.Cells.Item(68,7).Value = 10
$ws2.Cells.Item(68,8).Value = 3

$ws2Src = $ws2.Range("A65:H65")
$ws2Dst66 = $ws2.Range("A66:H66")
Restore-RowStyle $ws2Src $ws2Dst66
$ws2Dst67 = $ws2.Range("A67:H67")
Restore-RowStyle $ws2Src $ws2Dst67
$ws2Dst68 = $ws2.Range("A68:H68")
Restore-RowStyle $ws2Src $ws2Dst68

# ===========================================================================
# Sheet 3: "Positive Keywords"
# ===========================================================================
$ws3 = $wb.Worksheets.Item("Positive Keywords")
$ws3Src = $ws3.Range("A103:G103")

$posRows = @(
    @(104, "Venkata Rami Reddy",      6, 6, "100.0%", "0.0%", "0.0%", "0.875"),
    @(105, "AP Government Employees", 6, 6, "100.0%", "0.0%", "0.0%", "0.855"),
    @(106, "PRC fitment",             6, 6, "100.0%", "0.0%", "0.0%", "0.875"),
    @(107, "Dearness Allowance",      6, 6, "100.0%", "0.0%", "0.0%", "0.860"),
    @(108, "political criticism",     6, 6, "100.0%", "0.0%", "0.0%", "0.855"),
    @(109, "employee welfare",        6, 6, "100.0%", "0.0%", "0.0%", "0.840")
)
foreach ($r in $posRows) {
    $row = $r[0]
    Set-TextCell $ws3 $row 1 $r[1]
    $ws3.Cells.Item($row,2).Value = $r[2]
    $ws3.Cells.Item($row,3).Value = $r[3]
    Set-TextCell $ws3 $row 4 $r[4]
    Set-TextCell $ws3 $row 5 $r[5]
    Set-TextCell $ws3 $row 6 $r[6]
    Set-TextCell $ws3 $row 7 $r[7]
    $dst = $ws3.Range(("A" + $row + ":G" + $row))
    Restore-RowStyle $ws3Src $dst
}

# ===========================================================================
# Sheet 4: "Negative Keywords"
# ===========================================================================
$ws4 = $wb.Worksheets.Item("Negative Keywords")
$ws4Src = $ws4.Range("A167:G167")

$negRows = @(
    @(168, "Jagan Letter",           15, 15, "100.0%", "0.0%", "0.0%", "0.863"),
    @(169, "AP Assembly",            15, 15, "100.0%", "0.0%", "0.0%", "0.863"),
    @(170, "Political Crisis",       14, 14, "100.0%", "0.0%", "0.0%", "0.850"),
    @(171, "Naidu vs Jagan",         13, 13, "100.0%", "0.0%", "0.0%", "0.860"),
    @(172, "Opposition Rights",      14, 14, "100.0%", "0.0%", "0.0%", "0.843"),
    @(173, "YSRC Boycott",           13, 13, "100.0%", "0.0%", "0.0%", "0.860"),
    @(174, "AP Assembly Session",    15, 15, "100.0%", "0.0%", "0.0%", "0.850"),
    @(175, "Atchannaidu",             7,  7, "100.0%", "0.0%", "0.0%", "0.855"),
    @(176, "onion farmers",           7,  7, "100.0%", "0.0%", "0.0%", "0.845"),
    @(177, "political news",          7,  7, "100.0%", "0.0%", "0.0%", "0.855"),
    @(178, "Telugu states",           7,  7, "100.0%", "0.0%", "0.0%", "0.850"),
    @(179, "agricultural policies",   7,  7, "100.0%", "0.0%", "0.0%", "0.845")
)
foreach ($r in $negRows) {
    $row = $r[0]
    Set-TextCell $ws4 $row 1 $r[1]
    $ws4.Cells.Item($row,2).Value = $r[2]
    $ws4.Cells.Item($row,3).Value = $r[3]
    Set-TextCell $ws4 $row 4 $r[4]
    Set-TextCell $ws4 $row 5 $r[5]
    Set-TextCell $ws4 $row 6 $r[6]
    Set-TextCell $ws4 $row 7 $r[7]
    $dst = $ws4.Range(("A" + $row + ":G" + $row))
    Restore-RowStyle $ws4Src $dst
}

# ===========================================================================
# Sheet 7: "Negative Issues"
# ===========================================================================
$ws7 = $wb.Worksheets.Item("Negative Issues")
$ws7Src = $ws7.Range("A126:E126")

$issueRows = @(
    @(127, "Government Employees' PRC Fitment and DA Issues", "Political Governance", "0.900"),
    @(128, "Unfulfilled Arrears and Reduction in Interim Relief", "Economic Development", "0.850"),
    @(129, "Protests Against Government Policy on Employee Financial Benefits", "Political Governance", "0.800"),
    @(130, "Allegations Against Current Government's Handling of Employee Welfare", "Social Welfare Schemes", "0.750"),
    @(131, "Opposition Rights and Legislative Procedures", "Political Governance", "0.900"),
    @(132, "Corruption Allegations Against Leaders", "Corruption & Transparency", "0.850"),
    @(133, "Ineffective Assembly Discussions", "Political Governance", "0.800"),
    @(134, "Prioritization of Political Drama over Public Issues", "Political Governance", "0.750"),
    @(135, "Criticism of YS Jagan's leadership regarding agricultural support, particularly concerning fertilizer availability.", "Agricultural Support", "0.900"),
    @(136, "Comments suggesting a failure in agricultural governance under the current leadership.", "Political Governance", "0.850"),
    @(137, "Public dissatisfaction with the effectiveness of the current government's initiatives impacting farmers.", "Economic Development", "0.800")
)
foreach ($r in $issueRows) {
    $row = $r[0]
    Set-TextCell $ws7 $row 1 $r[1]
    Set-TextCell $ws7 $row 2 $r[2]
    Set-TextCell $ws7 $row 3 $r[3]
    $ws7.Cells.Item($row,4).Value = 1
    $ws7.Cells.Item($row,5).Value = 1
    $dst = $ws7.Range(("A" + $row + ":E" + $row))
    Restore-RowStyle $ws7Src $dst
}

# ===========================================================================
# Sheet 9: "Narrative Summary"
# ===========================================================================
$ws9 = $wb.Worksheets.Item("Narrative Summary")
$summaryText = "=== Updated 2025-10-16 21:26:27 ===`n`n" + `
    "AI-Generated Keyword Analysis Summary (Cumulative)`n`n" + `
    "The sentiment analysis of 70 YouTube videos related to Andhra Pradesh's political landscape reveals a strikingly negative public perception, with a complete absence of positive sentiment among the processed 52 comments. This overwhelming negativity, quantified at 100%, indicates a profound discontent among viewers regarding the political situation in the state.`n`n" + `
    "Key findings underscore critical themes shaping public discourse. First, the dominance of keywords such as 'Atchannaidu,' 'YS Jagan,' and 'onion farmers' suggests a significant concern regarding agricultural policies and the leadership of both the current and former chief ministers. The explicit mention of 'onion farmers' reflects specific grievances tied to agricultural distress, highlighting a potentially volatile issue that could impact electoral prospects.`n`n" + `
    "Furthermore, the polarizing figures of 'Chandrababu Naidu' and 'Pawan Kalyan' indicate a contentious political climate, where public figures elicit strong emotional responses. The absence of positive keywords reinforces the notion that political narratives are currently dominated by criticism rather than support, suggesting that both major parties may struggle to cultivate favorable public sentiment.`n`n" + `
    "This analysis reveals an urgent need for political actors in Andhra Pradesh to address the prevailing discontent, particularly in agricultural policy and leadership accountability. Engaging with these critical issues may be essential for any party seeking to shift public sentiment and regain trust in a highly charged political environment.`n`n" + `
    "This cumulative summary was generated by AI based on aggregated sentiment data from 70 videos and 52 comments."
$ws9.Cells.Item(1,1).Formula = $summaryText

# ===========================================================================
# Sheet 10: "Political Issues Narrative"
# ===========================================================================
$ws10 = $wb.Worksheets.Item("Political Issues Narrative")
$apostrophe = [char]0x2019
$politicalText = "=== Updated 2025-10-16 21:26:33 ===`n`n" + `
    "AI-Generated Political Issues Analysis (Cumulative)`n`n" + `
    "The analysis of 70 YouTube videos focusing on Andhra Pradesh politics reveals significant public discontent, particularly surrounding three key issues: agricultural support, political governance, and economic development. These issues are categorized under negative sentiment, indicating a prevailing dissatisfaction with the current administration led by YS Jagan Mohan Reddy.`n`n" + `
    "One of the most prominent concerns is the criticism of Jagan's leadership regarding agricultural support, specifically the availability of fertilizers. This issue underscores a broader narrative of perceived governmental failure in addressing the needs of farmers, a crucial demographic in Andhra Pradesh's economy. Comments reflecting a failure in agricultural governance further exacerbate public sentiment, highlighting a disconnect between the government" + $apostrophe + "s initiatives and the on-ground realities faced by the agricultural community.`n`n" + `
    "Moreover, public dissatisfaction with the overall effectiveness of government initiatives impacting farmers suggests a critical examination of economic development policies. This sentiment not only reveals the challenges in policy implementation but also signals a potential risk for the ruling party, as agricultural distress can have far-reaching social and electoral implications.`n`n" + `
    "In the context of Andhra Pradesh's political landscape, these findings illuminate the urgent need for the government to address these contentious issues. The lack of positive sentiment surrounding these themes indicates a pressing demand for accountability and effective governance. As discourse shifts towards these priorities, it is essential for policymakers to recalibrate their strategies, ensuring that they resonate with the public's concerns, particularly in agricultural and economic spheres.`n`n" + `
    "This cumulative analysis was generated by AI based on identified political issues from 70 videos."
$ws10.Cells.Item(1,1).Formula = $politicalText
